$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Duplicate the sheet (this is how the sheetId counter advances to 2,
# matching the saved file's internal bookkeeping) and drop the original,
# keeping the same name/position/relationship id.
$null = $ws.Copy($null, $ws)
$null = $wb.Worksheets.Item(1).Delete()
$wb.Worksheets.Item(1).Name = "Dummy"

$ws = $wb.Worksheets.Item(1)
$null = $ws.Activate()

# Add the new row of data.
$ws.Range("A16").Value = 1
$ws.Range("B16").Value = 2

# Move the selection the way it was left in the saved workbook.
$null = $ws.Range("A17").Select()
